# Updated cryptos list values (Price / Volume(1h)) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    # Force the literal text (avoids Excel auto-coercing numeric-looking
    # strings like "1.003" or "0.000008680" into numbers, which would
    # drop significant trailing zeros / switch to scientific notation).
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    # Restore the default (unstyled) formatting so only the value changes.
    $rng.ClearFormats()
}

Set-TextValue "D2" '27.202.76'
Set-TextValue "D3" '1.905.30'
Set-TextValue "E4" '  +0.30%  '
Set-TextValue "D5" '307.77'
Set-TextValue "E5" '  +0.59%  '
Set-TextValue "D6" '1.003'
Set-TextValue "E6" '  +0.27%  '
Set-TextValue "D7" '0.5265'
Set-TextValue "E7" '  +0.50%  '
Set-TextValue "D9" '0.07303'
Set-TextValue "D10" '21.56'
Set-TextValue "D11" '0.9052'
Set-TextValue "E11" '  +0.50%  '
Set-TextValue "D12" '0.08089'
Set-TextValue "E12" '  -4.71%  '
Set-TextValue "D13" '95.92'
Set-TextValue "D14" '5.374'
Set-TextValue "E14" '  +1.66%  '
Set-TextValue "D15" '1.804.92'
Set-TextValue "E15" '  -4.63%  '
Set-TextValue "E16" '  +0.26%  '
Set-TextValue "D17" '0.000008680'
Set-TextValue "E17" '  +0.65%  '
Set-TextValue "E18" '  +1.18%  '
Set-TextValue "E19" '  +0.25%  '
Set-TextValue "D20" '27.240.92'
Set-TextValue "E20" '  +0.12%  '
Set-TextValue "D21" '5.124'
Set-TextValue "E21" '  +1.10%  '
Set-TextValue "D22" '10.83'
Set-TextValue "E22" '  +2.11%  '
Set-TextValue "D23" '6.495'
Set-TextValue "E23" '  +1.11%  '
Set-TextValue "D24" '2.340'
Set-TextValue "E24" '  +2.55%  '
Set-TextValue "D25" '150.14'
Set-TextValue "E25" '  +1.94%  '
Set-TextValue "D26" '18.25'
Set-TextValue "E26" '  +0.41%  '
Set-TextValue "D27" '1.742'
Set-TextValue "E27" '  -0.45%  '
Set-TextValue "D28" '116.89'
Set-TextValue "E28" '  +1.73%  '
Set-TextValue "D29" '4.845'
Set-TextValue "E29" '  +0.62%  '
Set-TextValue "D30" '4.873'
Set-TextValue "E30" '  -0.37%  '
Set-TextValue "D31" '0.09228'
Set-TextValue "E31" '  -0.39%  '
Set-TextValue "D32" '0.8121'
Set-TextValue "E32" '  +0.53%  '
Set-TextValue "D33" '0.05068'
Set-TextValue "E33" '  +0.10%  '
Set-TextValue "D34" '1.230'
Set-TextValue "E34" '  -0.56%  '
Set-TextValue "D35" '2.984'
Set-TextValue "E35" '  +1.33%  '
Set-TextValue "D36" '3.362'
Set-TextValue "E36" '  -1.95%  '
Set-TextValue "D37" '2.702'
Set-TextValue "E37" '  +3.19%  '
Set-TextValue "D38" '0.5728'
Set-TextValue "E38" '  +0.03%  '
Set-TextValue "D39" '0.01995'
Set-TextValue "E39" '  +0.27%  '
Set-TextValue "D40" '1.087'
Set-TextValue "E40" '  +1.08%  '
Set-TextValue "D41" '8.981'
Set-TextValue "E41" '  -0.43%  '
Set-TextValue "D43" '116.89'
Set-TextValue "E43" '  +0.27%  '
Set-TextValue "E44" '  +0.53%  '
Set-TextValue "D45" '0.4935'
Set-TextValue "E45" '  +1.51%  '
Set-TextValue "E46" '  +0.26%  '
Set-TextValue "E47" '  -0.15%  '
Set-TextValue "D48" '1.637'
Set-TextValue "E48" '  +1.45%  '
Set-TextValue "D49" '38.54'
Set-TextValue "E49" '  +2.89%  '
Set-TextValue "D50" '64.32'
Set-TextValue "E50" '  +0.60%  '
Set-TextValue "D51" '0.05966'
Set-TextValue "E51" '  +0.28%  '
